$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AIKO", "ARON Unicorn Toilet", "10", "106.47", "1064.7"),
    @("Telerik Academy", "IKEA Dendroid Chair", "52", "73.8", "3837.6"),
    @("Boyana Film Studios", "ARON Medusa Mirror", "24", "52.74", "1265.76"),
    @("San Benedetto", "IKEA Medusa Mirror", "8", "36.95", "295.6"),
    @("Mladost Estate", "ARON Dendroid Chair", "4", "213.85", "855.4"),
    @("Sunset Security", "ARON Vampire Lamp", "14", "86.16", "1206.24"),
    @("West Bank", "ARON Dendroid Chair", "20", "213.85", "4277"),
    @("BILLA", "ARON Zombie Bed", "17", "193.05", "3281.85"),
    @("Null Industries", "IKEA Dendroid Chair", "21", "73.8", "1549.8"),
    @("VS Incorporated", "ARON Vampire Lamp", "7", "86.16", "603.12")
)

$r = 122
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
    $ws.Range("A" + $r + ":E" + $r).Style = "Normal"
    $r = $r + 1
}
